# Update 12-Jul-2021, end of day update.
# Replaces the week of 5-Jul..9-Jul-2021 entries on "Sheet1" (rows 3-35) with
# a single day's entries for 12-Jul-2021 (rows 3-8), and clears out the now
# unused trailing rows (9-35) so only the running "Saldo" formula remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Opening balance for the new period -----------------------------------
$ws.Range("E2").Value = 271225

# --- Row 3: Wages Expense (date changes 5-Jul -> 12-Jul-2021) -------------
$ws.Range("A3").Value = 44389
$ws.Range("D3").Formula = "=60000"

# --- Row 4: was "A/R" -> now "BELI - beras" --------------------------------
$ws.Range("B4").Value = "BELI - beras"
$ws.Range("C4").Clear()
$ws.Range("D4").Value = 290000

# --- Row 5: was "TRANSFER BCA" -> now "FREIGHT OUT" ------------------------
$ws.Range("B5").Value = "FREIGHT OUT"
$ws.Range("D5").Formula = "=144000"

# --- Row 6: was "FREIGHT - OUT" -> now "TRANSFER BCA" ----------------------
$ws.Range("B6").Value = "TRANSFER BCA"
$ws.Range("D6").Formula = "=406000+280000+4600000+5000000+3000000+41430000+3490000+580000+6560000"

# --- Row 7: was "PLN - Astar 214" -> now "A/R" ------------------------------
$ws.Range("B7").Value = "A/R"
$ws.Range("D7").Clear()
$ws.Range("C7").Formula = "=5000000+875000+41430000+6560000"
# F7 keeps its pre-existing (empty, styled) cell untouched.

# --- Row 8: was "SALES - cash/retail" -> now "A/P" --------------------------
$ws.Range("B8").Value = "A/P"
$ws.Range("C8").Clear()
$ws.Range("D8").Formula = "=1051200"

# --- Rows 9-16: drop the rest of the old week's entries (keep col E) -------
for ($r = 9; $r -le 16; $r++) {
    $ws.Range("A" + $r + ":D" + $r).Clear()
}

# --- Rows 17-35: drop the rest of the old week's entries (keep col E) ------
for ($r = 17; $r -le 35; $r++) {
    $ws.Range("A" + $r + ":D" + $r).Clear()
}

# --- Update the frozen-pane view to show the top of the sheet, with B9 ----
# --- selected (matches the end-of-day editing position). -------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B9").Select()
